$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 3939.3635
$ws.Cells.Item(43, 9).Value = 3638.8333
$ws.Cells.Item(43, 11).Value = 3638.8333
$ws.Cells.Item(43, 13).Value = -3569.8333
$ws.Cells.Item(62, 8).Value = 4309.3335
$ws.Cells.Item(62, 9).Value = 3464.1667
$ws.Cells.Item(62, 10).Value = 5999.6665
$ws.Cells.Item(62, 11).Value = 3464.1667
$ws.Cells.Item(62, 12).Value = 5999.6665
$ws.Cells.Item(62, 13).Value = -2840.1667
$ws.Cells.Item(62, 14).Value = -7247.6665
$ws.Cells.Item(65, 8).Value = 4309.3335
$ws.Cells.Item(65, 9).Value = 3464.1667
$ws.Cells.Item(65, 10).Value = 5999.6665
$ws.Cells.Item(65, 11).Value = 17320.8335
$ws.Cells.Item(65, 12).Value = 29998.3325
$ws.Cells.Item(65, 13).Value = -14200.8335
$ws.Cells.Item(65, 14).Value = -36238.3325
$ws.Cells.Item(92, 8).Value = 25001070
$ws.Cells.Item(92, 9).Value = 29412830
$ws.Cells.Item(92, 10).Value = 1097.3334
$ws.Cells.Item(92, 11).Value = 29412830
$ws.Cells.Item(92, 12).Value = 1097.3334
$ws.Cells.Item(92, 13).Value = -29411582
$ws.Cells.Item(92, 14).Value = -3593.3334
$ws.Cells.Item(96, 8).Value = 2191.4443
$ws.Cells.Item(96, 9).Value = 3158.4
$ws.Cells.Item(96, 10).Value = 982.75
$ws.Cells.Item(96, 11).Value = 9475.200000000001
$ws.Cells.Item(96, 12).Value = 2948.25
$ws.Cells.Item(96, 13).Value = -8102.200000000001
$ws.Cells.Item(96, 14).Value = -5694.25
$ws.Cells.Item(113, 8).Value = 21885.3
$ws.Cells.Item(113, 9).Value = 26669.25
$ws.Cells.Item(113, 10).Value = 2749.5
$ws.Cells.Item(113, 11).Value = 26669.25
$ws.Cells.Item(113, 12).Value = 2749.5
$ws.Cells.Item(113, 13).Value = -23415.25
$ws.Cells.Item(113, 14).Value = -9257.5
$ws.Cells.Item(132, 8).Value = 1581.5319
$ws.Cells.Item(132, 9).Value = 1522.8
$ws.Cells.Item(132, 11).Value = 4568.4
$ws.Cells.Item(132, 13).Value = -2038.4
$ws.Cells.Item(135, 8).Value = 3158.3076
$ws.Cells.Item(135, 9).Value = 1482.5
$ws.Cells.Item(135, 10).Value = 5839.6
$ws.Cells.Item(135, 11).Value = 13342.5
$ws.Cells.Item(135, 12).Value = 52556.4
$ws.Cells.Item(135, 13).Value = -10807.5
$ws.Cells.Item(135, 14).Value = -57626.4
$ws.Cells.Item(138, 8).Value = 2459.0808
$ws.Cells.Item(138, 9).Value = 1213.8889
$ws.Cells.Item(138, 10).Value = 2735.79
$ws.Cells.Item(138, 11).Value = 3641.6667
$ws.Cells.Item(138, 12).Value = 8207.369999999999
$ws.Cells.Item(138, 13).Value = 1498.3333
$ws.Cells.Item(138, 14).Value = -18487.37

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 30098050
$ws.Cells.Item(32, 9).Value = 31353464
$ws.Cells.Item(32, 11).Value = 31353464
$ws.Cells.Item(32, 13).Value = -31353177
$ws.Cells.Item(45, 8).Value = 3090.4
$ws.Cells.Item(45, 9).Value = 976
$ws.Cells.Item(45, 11).Value = 976
$ws.Cells.Item(45, 13).Value = -599
$ws.Cells.Item(74, 8).Value = 2279.6365
$ws.Cells.Item(74, 9).Value = 2230.9333
$ws.Cells.Item(74, 11).Value = 2230.9333
$ws.Cells.Item(74, 13).Value = -1356.9333
$ws.Cells.Item(77, 8).Value = 2279.6365
$ws.Cells.Item(77, 9).Value = 2230.9333
$ws.Cells.Item(77, 11).Value = 11154.6665
$ws.Cells.Item(77, 13).Value = -6786.666500000001
$ws.Cells.Item(97, 8).Value = 1436.375
$ws.Cells.Item(97, 9).Value = 944.75
$ws.Cells.Item(97, 10).Value = 2419.625
$ws.Cells.Item(97, 11).Value = 944.75
$ws.Cells.Item(97, 12).Value = 2419.625
$ws.Cells.Item(97, 13).Value = -448.75
$ws.Cells.Item(97, 14).Value = -3411.625
$ws.Cells.Item(132, 8).Value = 2287.096
$ws.Cells.Item(132, 9).Value = 2234.111
$ws.Cells.Item(132, 11).Value = 6702.333
$ws.Cells.Item(132, 13).Value = -4172.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3339.1
$ws.Cells.Item(20, 9).Value = 2876.6
$ws.Cells.Item(20, 10).Value = 3801.6
$ws.Cells.Item(20, 11).Value = 2876.6
$ws.Cells.Item(20, 12).Value = 3801.6
$ws.Cells.Item(20, 13).Value = -2629.6
$ws.Cells.Item(20, 14).Value = -4295.6
$ws.Cells.Item(86, 8).Value = 2221.2666
$ws.Cells.Item(86, 9).Value = 2612.5
$ws.Cells.Item(86, 11).Value = 2612.5
$ws.Cells.Item(86, 13).Value = -1489.5
$ws.Cells.Item(89, 8).Value = 2221.2666
$ws.Cells.Item(89, 9).Value = 2612.5
$ws.Cells.Item(89, 11).Value = 13062.5
$ws.Cells.Item(89, 13).Value = -7446.5
$ws.Cells.Item(134, 8).Value = 1625420.9
$ws.Cells.Item(134, 9).Value = 1881092.5
$ws.Cells.Item(134, 10).Value = 6166.5
$ws.Cells.Item(134, 11).Value = 5643277.5
$ws.Cells.Item(134, 12).Value = 18499.5
$ws.Cells.Item(134, 13).Value = -5640742.5
$ws.Cells.Item(134, 14).Value = -23569.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(23, 8).Value = 5500
$ws.Cells.Item(23, 9).Value = 5500
$ws.Cells.Item(23, 11).Value = 5500
$ws.Cells.Item(23, 13).Value = -5260
$ws.Cells.Item(27, 8).Value = 5500
$ws.Cells.Item(27, 9).Value = 5500
$ws.Cells.Item(27, 11).Value = 5500
$ws.Cells.Item(27, 13).Value = -5308
$ws.Cells.Item(31, 8).Value = 3708.4546
$ws.Cells.Item(31, 10).Value = 3935.4644
$ws.Cells.Item(31, 12).Value = 3935.4644
$ws.Cells.Item(31, 14).Value = -4525.4644
$ws.Cells.Item(34, 8).Value = 3708.4546
$ws.Cells.Item(34, 10).Value = 3935.4644
$ws.Cells.Item(34, 12).Value = 3935.4644
$ws.Cells.Item(34, 14).Value = -4339.4644
$ws.Cells.Item(107, 8).Value = 1649.5264
$ws.Cells.Item(107, 10).Value = 2265.5454
$ws.Cells.Item(107, 12).Value = 2265.5454
$ws.Cells.Item(107, 14).Value = -6105.5454
$ws.Cells.Item(122, 8).Value = 3758.3438
$ws.Cells.Item(122, 9).Value = 3327.5
$ws.Cells.Item(122, 10).Value = 4189.1875
$ws.Cells.Item(122, 11).Value = 9982.5
$ws.Cells.Item(122, 12).Value = 12567.5625
$ws.Cells.Item(122, 13).Value = -7532.5
$ws.Cells.Item(122, 14).Value = -17467.5625
$ws.Cells.Item(132, 8).Value = 3191.718
$ws.Cells.Item(132, 9).Value = 3235.125
$ws.Cells.Item(132, 11).Value = 9705.375
$ws.Cells.Item(132, 13).Value = -7175.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 89424110
$ws.Cells.Item(4, 9).Value = 57406108
$ws.Cells.Item(4, 11).Value = 172218324
$ws.Cells.Item(4, 13).Value = -172218212
$ws.Cells.Item(18, 8).Value = 3019.2666
$ws.Cells.Item(18, 9).Value = 1929
$ws.Cells.Item(18, 11).Value = 5787
$ws.Cells.Item(18, 13).Value = -5618
$ws.Cells.Item(38, 8).Value = 45.82353
$ws.Cells.Item(38, 10).Value = 72.40000000000001
$ws.Cells.Item(38, 12).Value = 217.2
$ws.Cells.Item(38, 14).Value = -911.2
$ws.Cells.Item(55, 8).Value = 3041.5
$ws.Cells.Item(55, 9).Value = 2200
$ws.Cells.Item(55, 10).Value = 3251.875
$ws.Cells.Item(55, 11).Value = 6600
$ws.Cells.Item(55, 12).Value = 9755.625
$ws.Cells.Item(55, 13).Value = -6423
$ws.Cells.Item(55, 14).Value = -10109.625
$ws.Cells.Item(113, 8).Value = 1807.125
$ws.Cells.Item(113, 10).Value = 1939.7142
$ws.Cells.Item(113, 12).Value = 5819.142599999999
$ws.Cells.Item(113, 14).Value = -10159.1426
$ws.Cells.Item(131, 8).Value = 1417.2858
$ws.Cells.Item(131, 10).Value = 1832
$ws.Cells.Item(131, 12).Value = 5496
$ws.Cells.Item(131, 14).Value = -15576
$ws.Cells.Item(139, 8).Value = 6616.1113
$ws.Cells.Item(139, 10).Value = 8581.091
$ws.Cells.Item(139, 12).Value = 25743.273
$ws.Cells.Item(139, 14).Value = -36023.273

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 108994
$ws.Cells.Item(39, 10).Value = 108994
$ws.Cells.Item(39, 12).Value = 108994
$ws.Cells.Item(39, 14).Value = -110058
$ws.Cells.Item(132, 8).Value = 2193.125
$ws.Cells.Item(132, 9).Value = 1844.1025
$ws.Cells.Item(132, 10).Value = 3705.5557
$ws.Cells.Item(132, 11).Value = 5532.3075
$ws.Cells.Item(132, 12).Value = 11116.6671
$ws.Cells.Item(132, 13).Value = -3002.3075
$ws.Cells.Item(132, 14).Value = -16176.6671
$ws.Cells.Item(136, 8).Value = 66586.17999999999
$ws.Cells.Item(136, 10).Value = 66586.17999999999
$ws.Cells.Item(136, 12).Value = 199758.54
$ws.Cells.Item(136, 14).Value = -204858.54

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 7320.3438
$ws.Cells.Item(46, 9).Value = 3250
$ws.Cells.Item(46, 10).Value = 7591.7
$ws.Cells.Item(46, 11).Value = 3250
$ws.Cells.Item(46, 12).Value = 7591.7
$ws.Cells.Item(46, 13).Value = -3062
$ws.Cells.Item(46, 14).Value = -7967.7
$ws.Cells.Item(122, 8).Value = 26402.545
$ws.Cells.Item(122, 9).Value = 27429.375
$ws.Cells.Item(122, 11).Value = 82288.125
$ws.Cells.Item(122, 13).Value = -79838.125
$ws.Cells.Item(132, 8).Value = 4401.857
$ws.Cells.Item(132, 9).Value = 4091.375
$ws.Cells.Item(132, 10).Value = 4986.294
$ws.Cells.Item(132, 11).Value = 12274.125
$ws.Cells.Item(132, 12).Value = 14958.882
$ws.Cells.Item(132, 13).Value = -9744.125
$ws.Cells.Item(132, 14).Value = -20018.882
$ws.Cells.Item(136, 8).Value = 13736.846
$ws.Cells.Item(136, 9).Value = 9823.375
$ws.Cells.Item(136, 10).Value = 19998.4
$ws.Cells.Item(136, 11).Value = 29470.125
$ws.Cells.Item(136, 12).Value = 59995.2
$ws.Cells.Item(136, 13).Value = -26920.125
$ws.Cells.Item(136, 14).Value = -65095.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(58, 8).Value = 9691
$ws.Cells.Item(58, 9).Value = 8613.75
$ws.Cells.Item(58, 11).Value = 8613.75
$ws.Cells.Item(58, 13).Value = -8305.75
$ws.Cells.Item(122, 8).Value = 7414.5713
$ws.Cells.Item(122, 10).Value = 2698
$ws.Cells.Item(122, 12).Value = 8094
$ws.Cells.Item(122, 14).Value = -12994
$ws.Cells.Item(132, 8).Value = 1660.25
$ws.Cells.Item(132, 9).Value = 1541.3939
$ws.Cells.Item(132, 10).Value = 2967.6667
$ws.Cells.Item(132, 11).Value = 4624.1817
$ws.Cells.Item(132, 12).Value = 8903.000100000001
$ws.Cells.Item(132, 13).Value = -2094.1817
$ws.Cells.Item(132, 14).Value = -13963.0001
$ws.Cells.Item(136, 8).Value = 8376506.5
$ws.Cells.Item(136, 10).Value = 33502666
$ws.Cells.Item(136, 12).Value = 100507998
$ws.Cells.Item(136, 14).Value = -100513098
